$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.755.78"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.800.13"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.65"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.23"
$ws.Range("E6").Value = "  +0.97%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  +0.48%  "
$ws.Range("E9").Value = "  +1.60%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.29"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.442.96"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.800.52"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.56"
$ws.Range("E16").Value = "  +3.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.810.28"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "461.07"
$ws.Range("E20").Value = "  +0.81%  "
$ws.Range("E21").Value = "  -2.80%  "
$ws.Range("E22").Value = "  -0.20%  "
$ws.Range("E23").Value = "  +1.39%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.07"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.946.95"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.76"
$ws.Range("E30").Value = "  -0.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.37"
$ws.Range("E31").Value = "  +2.88%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.52"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.741.65"
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("E38").Value = "  +2.83%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("E41").Value = "  +1.22%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.11"
$ws.Range("E44").Value = "  +2.53%  "
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.65"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.33"
$ws.Range("E47").Value = "  +0.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.30"
$ws.Range("E48").Value = "  +7.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "147.64"
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  +9.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "394.51"
$ws.Range("E51").Value = "  +0.71%  "
